# Generate Report for Handback
# Populates the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns on the zh-cn and de-de detail sheets now that the handback has completed, updates
# the Status text, and widens a few columns to fit the new (longer) content.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c0fd6f980678a7f9d0450f60b918c418fc919434/e2e"

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("J2").Value = "39f55124-af5e-4722-8f16-5825e37e1315.aee9819775a902ca9c4058e10af767c457a615a4.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-19 06:24:54"

$wsZh.Range("J3").Value = "b25f734d-6242-45c0-be2b-ea661aa09933.6a947351c110a4a4f6ca4730c5917fd23d76c810.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-19 06:24:54"

# Rebuild the hyperlinks collection so the new "Latest Target File" links land in the
# right order (A2, I2, A3, I3) with a fresh set of relationship ids.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$repoBase/39f55124-af5e-4722-8f16-5825e37e1315.md", "", "", "39f55124-af5e-4722-8f16-5825e37e1315.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "$repoBase/39f55124-af5e-4722-8f16-5825e37e1315.md", "", "", "39f55124-af5e-4722-8f16-5825e37e1315.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "$repoBase/b25f734d-6242-45c0-be2b-ea661aa09933.md", "", "", "b25f734d-6242-45c0-be2b-ea661aa09933.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "$repoBase/b25f734d-6242-45c0-be2b-ea661aa09933.md", "", "", "b25f734d-6242-45c0-be2b-ea661aa09933.md") | Out-Null

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("J2").Value = "39f55124-af5e-4722-8f16-5825e37e1315.aee9819775a902ca9c4058e10af767c457a615a4.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-19 06:25:04"

$wsDe.Range("J3").Value = "b25f734d-6242-45c0-be2b-ea661aa09933.6a947351c110a4a4f6ca4730c5917fd23d76c810.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-19 06:25:04"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$repoBase/39f55124-af5e-4722-8f16-5825e37e1315.md", "", "", "39f55124-af5e-4722-8f16-5825e37e1315.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "$repoBase/39f55124-af5e-4722-8f16-5825e37e1315.md", "", "", "39f55124-af5e-4722-8f16-5825e37e1315.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$repoBase/b25f734d-6242-45c0-be2b-ea661aa09933.md", "", "", "b25f734d-6242-45c0-be2b-ea661aa09933.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "$repoBase/b25f734d-6242-45c0-be2b-ea661aa09933.md", "", "", "b25f734d-6242-45c0-be2b-ea661aa09933.md") | Out-Null

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# Overview sheet - widen the zh-cn / de-de status columns (E, F) to match
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668
